$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> (predidx string, pred_name string)
$updates = @{
    3   = @("[1, 0, 0, 1, 0, 0, 0]", "['Normal', 'ParamViolation']")
    11  = @("[1, 0, 1, 0, 1, 0, 0]", "['Normal', 'HardwareFault', 'RegulationViolation']")
    12  = @("[1, 0, 1, 0, 0, 0, 0]", "['Normal', 'HardwareFault']")
    15  = @("[0, 0, 0, 0, 0, 0, 0]", "[]")
    22  = @("[1, 0, 1, 0, 0, 0, 0]", "['Normal', 'HardwareFault']")
    24  = @("[0, 0, 0, 0, 0, 0, 1]", "['SoftwareFault']")
    26  = @("[0, 0, 0, 0, 0, 0, 0]", "[]")
    29  = @("[1, 0, 0, 0, 0, 0, 1]", "['Normal', 'SoftwareFault']")
    31  = @("[1, 0, 0, 0, 0, 0, 1]", "['Normal', 'SoftwareFault']")
    36  = @("[1, 1, 0, 0, 0, 0, 0]", "['Normal', 'SurroundingEnvironment']")
    38  = @("[1, 0, 0, 0, 0, 0, 0]", "['Normal']")
    39  = @("[1, 0, 1, 0, 0, 0, 0]", "['Normal', 'HardwareFault']")
    61  = @("[0, 0, 0, 0, 0, 0, 1]", "['SoftwareFault']")
    68  = @("[1, 0, 0, 0, 0, 0, 0]", "['Normal']")
    70  = @("[1, 1, 0, 0, 0, 1, 0]", "['Normal', 'SurroundingEnvironment', 'CommunicationIssue']")
    71  = @("[1, 0, 0, 0, 0, 0, 0]", "['Normal']")
    73  = @("[1, 0, 0, 0, 0, 0, 1]", "['Normal', 'SoftwareFault']")
    80  = @("[1, 0, 1, 0, 0, 0, 0]", "['Normal', 'HardwareFault']")
    81  = @("[1, 0, 1, 0, 0, 0, 0]", "['Normal', 'HardwareFault']")
    84  = @("[1, 0, 0, 0, 0, 1, 0]", "['Normal', 'CommunicationIssue']")
    92  = @("[1, 0, 1, 0, 0, 0, 1]", "['Normal', 'HardwareFault', 'SoftwareFault']")
    97  = @("[0, 0, 0, 0, 0, 0, 0]", "[]")
    109 = @("[1, 1, 0, 0, 0, 0, 0]", "['Normal', 'SurroundingEnvironment']")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}
